$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Ghi chú" (note) column - header + per-member role
$ws.Range("E2").Value = "Ghi chú"
$ws.Range("E3").Value = "Nhóm trưởng"
$ws.Range("E4").Value = "Thành viên"
$ws.Range("E5").Value = "Thành viên"
$ws.Range("E6").Value = "Thành viên"
$ws.Range("E7").Value = "Thành viên"

# Extra reference lines below the table
$ws.Range("B9").Value = "Đường dẫn thùng chứa:http://quanlythuvien5n.googlecode.com/svn/trunk/"
$ws.Range("B10").Value = "Đường dẫn đến project: http://code.google.com/p/quanlythuvien5n/"

# Re-fit column widths for the updated content
$ws.Columns("A").ColumnWidth = 4
$ws.Columns("B").ColumnWidth = 22.714285714285715
$ws.Columns("C").ColumnWidth = 9.571428571428571
$ws.Columns("D").ColumnWidth = 29.428571428571427
$ws.Columns("E").ColumnWidth = 14.428571428571429

# Move the active selection
$ws.Range("B14").Select() | Out-Null
